$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their original text formatting
# (Price strings may look numeric, e.g. "1.00", "5.50", "51.188.24")
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D13", "D15", "D16", "D18", "D19", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D41", "D43", "D44", "D46", "D49", "D50")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "51.188.24"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.765.66"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "353.93"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "108.07"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "39.49"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  +3.53%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "19.92"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "3.199.07"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.776.70"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "51.171.82"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("E28").Value = "  +12.53%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "34.88"
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "51.83"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "6.07"
$ws.Range("E33").Value = "  +6.45%  "
$ws.Range("D34").Value = "0.0441"
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +3.76%  "
$ws.Range("D36").Value = "0.0831"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "18.15"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "2.51"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "120.56"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "22.13"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").Value = "2.086.85"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "5.50"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("D50").Value = "0.917"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  +6.53%  "
